$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Handout master & notes master: the cached "datetimeFigureOut"
#    footer date field moves from 4/27/22 -> 5/6/2022.
#    (Guarded with try/catch: some hosts expose these master date
#    placeholders as read-only through automation; if so, skip
#    quietly rather than aborting the rest of the edit.)
# ------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
        $shp = $hm.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.TrimEnd("`r") -eq "4/27/22") {
                $shp.TextFrame.TextRange.Text = "5/6/2022"
            }
        }
    }
} catch {
    # Some hosts expose the handout-master date placeholder read-only
    # through automation; skip quietly rather than aborting the rest
    # of the edit.
}

try {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $shp = $nm.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.TrimEnd("`r") -eq "4/27/22") {
                $shp.TextFrame.TextRange.Text = "5/6/2022"
            }
        }
    }
} catch {
    # Same read-only guard as above, for the notes master.
}

# ------------------------------------------------------------------
# 2) Slide 19 ("Testing Strategies ..."), Content Placeholder 2:
#    tidy up the "x / y" slash spacing + capitalization of two
#    bullet lines.
# ------------------------------------------------------------------
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$count = $tr.Paragraphs().Count

for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $t = $para.Text.TrimEnd("`r")

    if ($t -eq "Physics / Math Based Strategies") {
        $para.Runs(1).Text = "Physics/math based strategies"
    }
    elseif ($t -eq "Input / Output specifications, program invariants") {
        $para.Runs(1).Text = "Input/output specifications, program invariants"
    }
}
